$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Abril de 2020 a las 01:52"

# Update changed country rows: reordering per updated case counts,
# refreshed totals, and one newly appearing country (Timor Oriental, row 212).
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 310133
$ws.Cells.Item(4, 3).Value = 32972
$ws.Cells.Item(4, 4).Value = 14741
$ws.Cells.Item(4, 5).Value = 286950
$ws.Cells.Item(4, 6).Value = 8206
$ws.Cells.Item(4, 7).Value = 1038
$ws.Cells.Item(4, 8).Value = 8442

$ws.Cells.Item(43, 1).Value = "Panama"
$ws.Cells.Item(43, 2).Value = 1801
$ws.Cells.Item(43, 3).Value = 128
$ws.Cells.Item(43, 4).Value = 13
$ws.Cells.Item(43, 5).Value = 1742
$ws.Cells.Item(43, 6).Value = 50
$ws.Cells.Item(43, 7).Value = 5
$ws.Cells.Item(43, 8).Value = 46

$ws.Cells.Item(44, 1).Value = "Peru"
$ws.Cells.Item(44, 2).Value = 1746
$ws.Cells.Item(44, 3).Value = 151
$ws.Cells.Item(44, 4).Value = 914
$ws.Cells.Item(44, 5).Value = 759
$ws.Cells.Item(44, 6).Value = 88
$ws.Cells.Item(44, 7).Value = 12
$ws.Cells.Item(44, 8).Value = 73

$ws.Cells.Item(45, 1).Value = "Mexico"
$ws.Cells.Item(45, 2).Value = 1688
$ws.Cells.Item(45, 3).Value = 178
$ws.Cells.Item(45, 4).Value = 633
$ws.Cells.Item(45, 5).Value = 995
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(45, 7).Value = 10
$ws.Cells.Item(45, 8).Value = 60

$ws.Cells.Item(49, 1).Value = "Republica Dominicana"
$ws.Cells.Item(49, 2).Value = 1578
$ws.Cells.Item(49, 3).Value = 90
$ws.Cells.Item(49, 4).Value = 17
$ws.Cells.Item(49, 5).Value = 1484
$ws.Cells.Item(49, 6).Value = 147
$ws.Cells.Item(49, 7).Value = 9
$ws.Cells.Item(49, 8).Value = 77

$ws.Cells.Item(50, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(50, 2).Value = 1505
$ws.Cells.Item(50, 3).Value = 241
$ws.Cells.Item(50, 4).Value = 125
$ws.Cells.Item(50, 5).Value = 1370
$ws.Cells.Item(50, 6).Value = 2
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 10

$ws.Cells.Item(87, 1).Value = "Uruguay"
$ws.Cells.Item(87, 2).Value = 400
$ws.Cells.Item(87, 3).Value = 14
$ws.Cells.Item(87, 4).Value = 93
$ws.Cells.Item(87, 5).Value = 302
$ws.Cells.Item(87, 6).Value = 13
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 5

$ws.Cells.Item(112, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(112, 2).Value = 154
$ws.Cells.Item(112, 3).Value = 6
$ws.Cells.Item(112, 4).Value = 3
$ws.Cells.Item(112, 5).Value = 133
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 2
$ws.Cells.Item(112, 8).Value = 18

$ws.Cells.Item(144, 1).Value = "Puerto Rico"
$ws.Cells.Item(144, 2).Value = 39
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 1
$ws.Cells.Item(144, 5).Value = 36
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 2

$ws.Cells.Item(145, 1).Value = "Zambia"
$ws.Cells.Item(145, 2).Value = 39
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 36
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 1

$ws.Cells.Item(173, 1).Value = "Fiyi"
$ws.Cells.Item(173, 2).Value = 12
$ws.Cells.Item(173, 3).Value = 5
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 12
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(174, 1).Value = "Granada"
$ws.Cells.Item(174, 2).Value = 12
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 12
$ws.Cells.Item(174, 6).Value = 2
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

$ws.Cells.Item(179, 1).Value = "Surinam"
$ws.Cells.Item(179, 2).Value = 10
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 9
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 1

$ws.Cells.Item(180, 1).Value = "Mozambique"
$ws.Cells.Item(180, 2).Value = 10
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 1
$ws.Cells.Item(180, 5).Value = 9
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "Liberia"
$ws.Cells.Item(181, 2).Value = 10
$ws.Cells.Item(181, 3).Value = 3
$ws.Cells.Item(181, 4).Value = 3
$ws.Cells.Item(181, 5).Value = 6
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 1
$ws.Cells.Item(181, 8).Value = 1

$ws.Cells.Item(183, 1).Value = "Sudan"
$ws.Cells.Item(183, 2).Value = 10
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 2
$ws.Cells.Item(183, 5).Value = 6
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 2

$ws.Cells.Item(184, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(184, 2).Value = 9
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Suazilandia"
$ws.Cells.Item(185, 2).Value = 9
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

$ws.Cells.Item(186, 1).Value = "Republica del Chad"
$ws.Cells.Item(186, 2).Value = 9
$ws.Cells.Item(186, 3).Value = 1
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 9
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

$ws.Cells.Item(192, 1).Value = "Cabo Verde"
$ws.Cells.Item(192, 2).Value = 7
$ws.Cells.Item(192, 3).Value = 1
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 6
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 1

$ws.Cells.Item(193, 1).Value = "Somalia"
$ws.Cells.Item(193, 2).Value = 7
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 5).Value = 6
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(194, 2).Value = 7
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 1
$ws.Cells.Item(194, 5).Value = 6
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

$ws.Cells.Item(200, 1).Value = "Sierra Leona"
$ws.Cells.Item(200, 2).Value = 4
$ws.Cells.Item(200, 3).Value = 2
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 4
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0

$ws.Cells.Item(201, 1).Value = "Malaui"
$ws.Cells.Item(201, 2).Value = 4
$ws.Cells.Item(201, 3).Value = 1
$ws.Cells.Item(201, 4).Value = 0
$ws.Cells.Item(201, 5).Value = 4
$ws.Cells.Item(201, 6).Value = 0
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 0

$ws.Cells.Item(202, 1).Value = "Sahara Occidental"
$ws.Cells.Item(202, 2).Value = 4
$ws.Cells.Item(202, 3).Value = 4
$ws.Cells.Item(202, 4).Value = 0
$ws.Cells.Item(202, 5).Value = 4
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

$ws.Cells.Item(203, 1).Value = "Belice"
$ws.Cells.Item(203, 2).Value = 4
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 0
$ws.Cells.Item(203, 5).Value = 4
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

$ws.Cells.Item(204, 1).Value = "Botsuana"
$ws.Cells.Item(204, 2).Value = 4
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 0
$ws.Cells.Item(204, 5).Value = 3
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 1

$ws.Cells.Item(205, 1).Value = "Gambia"
$ws.Cells.Item(205, 2).Value = 4
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 2
$ws.Cells.Item(205, 5).Value = 1
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 1

$ws.Cells.Item(208, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(208, 2).Value = 3
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 0
$ws.Cells.Item(208, 5).Value = 3
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(209, 2).Value = 2
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 0
$ws.Cells.Item(209, 5).Value = 2
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Islas Malvinas"
$ws.Cells.Item(211, 2).Value = 1
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 0
$ws.Cells.Item(211, 5).Value = 1
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

$ws.Cells.Item(212, 1).Value = "Timor Oriental"
$ws.Cells.Item(212, 2).Value = 1
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 0
$ws.Cells.Item(212, 5).Value = 1
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 0

